# Daily_Motivations.xlsx - "Updating the daily scores"
# Applies:
#   1) For every existing "sleep" row (r = 2,5,8,...,65), flip the JKL (N)
#      and OS (O) score columns from TRUE to FALSE.
#   2) A handful of one-off score corrections on rows 66, 68, 69, 70
#      (Iron Man / F column, plus JKL+OS on row 68).
#   3) Append six new scored rows (71-76) for 2025-02-24 and 2025-02-25
#      (sleep / activity / weekly_activity for each date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Bulk JKL (N) / OS (O) correction on every "sleep" row 2..65 ---
for ($r = 2; $r -le 65; $r += 3) {
    $ws.Cells.Item($r, 14).Value = $false
    $ws.Cells.Item($r, 15).Value = $false
}

# --- 2) One-off corrections near the bottom of the existing data ---
$ws.Cells.Item(66, 6).Value = $true    # F66 (Iron Man / activity)  0 -> 1

$ws.Cells.Item(68, 6).Value = $false   # F68 (Iron Man / sleep)     1 -> 0
$ws.Cells.Item(68, 14).Value = $false  # N68 (JKL / sleep)          1 -> 0
$ws.Cells.Item(68, 15).Value = $false  # O68 (OS / sleep)           1 -> 0

$ws.Cells.Item(69, 6).Value = $true    # F69 (Iron Man / activity)  0 -> 1

$ws.Cells.Item(70, 6).Value = $true    # F70 (Iron Man / weekly)    0 -> 1

# --- 3) Append new rows 71-76 ---
function Set-ScoreRow {
    param($RowNum, $DateText, $MotivationTyp, $Scores)

    $dateCell = $ws.Cells.Item($RowNum, 1)
    # Force text so the date-like string isn't auto-converted to a serial
    # date number/format, then drop the temporary style so the cell is
    # left with no explicit format (matching the rest of the sheet).
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $DateText
    $dateCell.ClearFormats()

    $ws.Cells.Item($RowNum, 2).Value = $MotivationTyp

    for ($i = 0; $i -lt $Scores.Length; $i++) {
        $ws.Cells.Item($RowNum, 3 + $i).Value = [bool]$Scores[$i]
    }
}

# Columns C..O = Sportfinke, Taylor Atwood, Summerbody25, Iron Man, GurkenSalat,
#                yKing, StayStrong, WobblyWheel, ClearMind23, HealthQuest,
#                DeadliftCarror, JKL, OS
Set-ScoreRow 71 "2025-02-24" "sleep"           @(1,0,0,1,1,1,1,1,1,1,1,0,0)
Set-ScoreRow 72 "2025-02-24" "activity"        @(1,0,0,1,1,1,0,1,0,0,0,0,0)
Set-ScoreRow 73 "2025-02-24" "weekly_activity" @(1,0,0,1,1,0,1,1,0,0,0,1,0)
Set-ScoreRow 74 "2025-02-25" "sleep"           @(0,0,1,1,0,1,1,1,1,0,1,0,0)
Set-ScoreRow 75 "2025-02-25" "activity"        @(0,0,1,0,1,1,0,1,0,1,0,0,0)
Set-ScoreRow 76 "2025-02-25" "weekly_activity" @(1,0,1,1,1,0,1,1,0,0,0,1,0)
